$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 132
$lastCol = 8

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null -and $v -is [string] -and $v.Contains(";")) {
            $cell.Value2 = $v.Replace(";", ",")
        }
    }
}

$rng = $ws.Range("A1:H132")
$rng.AutoFilter() | Out-Null

$nm = $ws.Names.Add("_xlnm._FilterDatabase", $rng)
$nm.Visible = $false

$ws.Range("J4").Select()
